$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E4").Value = 6
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 13
$ws.Range("E4").Select()
